$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    ,@(0.2917716402565462, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 24.83232187738678)
    ,@(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    ,@(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    ,@(0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.069511820747223)
    ,@(0.04271373187048222, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.330410019770453)
    ,@(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    ,@(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    ,@(0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.069511820747223)
    ,@(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 7.143138311642302)
    ,@(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    ,@(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    ,@(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    ,@(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    ,@(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    ,@(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    ,@(1.455362044514542, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 3.009163075608874)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
    $ws.Cells.Item($row, 4).Value = $values[$i][2]
    $ws.Cells.Item($row, 5).Value = $values[$i][3]
    $ws.Cells.Item($row, 7).Value = $values[$i][4]
}

Write-Output "applied sval regen"
